$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.111.59"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.827.81"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.97"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  +7.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3743"
$ws.Range("E8").Value = "  +2.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07327"
$ws.Range("E9").Value = "  +1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8632"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.99"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "1.823.52"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.721"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.01"
$ws.Range("E14").Value = "  +5.87%  "
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07083"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008864"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "27.128.44"
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.205"
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.04"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.003"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.93"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.233"
$ws.Range("E26").Value = "  +5.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.52"
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.286"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.70"
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08919"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7652"
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("E33").Value = "  +5.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.482"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01974"
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05301"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5375"
$ws.Range("E39").Value = "  +7.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.199"
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("E42").Value = "  +2.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5209"
$ws.Range("E43").Value = "  +10.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.653"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.69"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.006"
$ws.Range("E46").Value = "  +11.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.21"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.689"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06450"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9252"
$ws.Range("E51").Value = "  +1.61%  "
